$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '24.514.01'
$ws.Range('E2').Value = '  -1.59%  '
$ws.Range('D3').Value = '1.668.77'
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '313.35'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.77%  '
$ws.Range('E6').Value = '  +0.05%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3896'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -3.36%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3925'
$ws.Range('D8').Style = 'Normal'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.004'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.03%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '51.61'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.86%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.398'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -4.85%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08612'
$ws.Range('D12').Style = 'Normal'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '24.98'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -3.81%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.272'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.82%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.00001314'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.69%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '7.711'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -4.11%  '
$ws.Range('D17').Value = '1.679.01'
$ws.Range('E17').Value = '  +1.51%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '93.08'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -3.51%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.07068'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.23%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '20.54'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.76%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '7.040'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.99%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.005'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.26%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '13.93'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.98%  '
$ws.Range('D24').Value = '24.508.30'
$ws.Range('E24').Value = '  -1.61%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.375'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.80%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '23.15'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.38%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.731'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -5.82%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '161.06'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.80%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.779'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -11.59%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '147.78'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.65%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.290'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.41%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.481'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +9.77%  '
$ws.Range('D33').Value = '1.891.00'
$ws.Range('E33').Value = '  +3.07%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08332'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -5.49%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.03014'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -6.31%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.944'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -5.03%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.2795'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.87%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.9778'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -4.35%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.09430'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.65%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.536'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +4.06%  '
$ws.Range('E41').Value = '  -5.17%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.7876'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -6.69%  '
$ws.Range('E43').Value = '  -4.44%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '16.42'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -7.88%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.7087'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -4.83%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.541'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -6.59%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.176'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.69%  '
$ws.Range('E48').Value = '  -0.52%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.08569'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.76%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.315'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -5.59%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '137.11'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.52%  '
